$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D: "success" styled like the existing B1/C1 headers
$ws.Range("D1").Value = "success"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# success flag (0/1, stored as text like the rest of the "list"/"count" data)
# per-row values taken from the target data
$successValues = @(
  0,0,0,0,0,0,0,0,0,0,0,0,0,0,
  1,1,0,1,1,0,1,1,1,1,1,1,1,
  0,0,1,1,0,0,0,0,1,1,1,1,1,
  1,0,1,1,0,1,1,0,0,0,0,0
)

$rng = $ws.Range("D2:D53")
$rng.NumberFormat = "@"

for ($i = 0; $i -lt $successValues.Length; $i++) {
    $row = 2 + $i
    $ws.Range("D$row").Value = [string]$successValues[$i]
}

$rng.ClearFormats()

Write-Output "done"
